$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content updates (per diff) ---

# E2 (Artificer / bonuses): "Arcane or Technology" -> "Arcane or Science"
$ws.Range("E2").Value = @'
\bonus{Imbue or Craft}{\twoCape}
\bonus{Arcane or Science}{\twoCape}
\bonus{Intelligence}{\twoCape}
\bonus{Logic}{\oneCape}
'@

# E7 (Responder / bonuses): "Kindness or Medicine"/"Medicine or Kindness" -> "Kindness or First Aid"/"First Aid or Kindness"
$ws.Range("E7").Value = @'
\bonus{Insight or Willpower}{\twoCape}
\bonus{Kindness or First Aid}{\twoCape}
\bonus{First Aid or Kindness}{\oneCape}
\bonus{Willpower or Insight}{\oneCape}
\bonus{Conviction}{\oneCape}
'@

# K7 (Responder / knowledgeDescription): "Medicine" -> "First Aid", drop "and identify symptoms and medical issues"
$ws.Range("K7").Value = @'
Whilst \imp{First Aid} allows you to identify and treats wounds and injuries, the advanced knowledge of the causes of diseases, infections, disorders and even magical curses is the domain of \imp{Pathology}. 

Though it might not help a broken ankle, a high \imp{Pathology} rating allows one to see the root cause of medical issues, and see the path needed to cure the malady. 
'@

# --- Row height adjustment for row 2 (144.75 -> 145.5) ---
$ws.Rows.Item(2).RowHeight = 145.5

# --- Selection / view update (topLeftCell A1->A5, active cell G3->K8) ---
$ws.Range("A5").Select()
$ws.Range("K8").Select()
